$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2
Set-TextValue "D2" "68.145.83"
Set-TextValue "E2" "  +1.30%  "

# Row 3
Set-TextValue "D3" "3.552.74"
Set-TextValue "E3" "  +1.87%  "

# Row 4
Set-TextValue "D4" "0.998"
Set-TextValue "E4" "  -0.22%  "

# Row 5
Set-TextValue "D5" "619.59"
Set-TextValue "E5" "  +2.38%  "

# Row 6
Set-TextValue "D6" "154.77"
Set-TextValue "E6" "  +4.41%  "

# Row 7
Set-TextValue "D7" "3.550.56"
Set-TextValue "E7" "  +1.81%  "

# Row 8
Set-TextValue "D8" "0.999"

# Row 9
Set-TextValue "E9" "  +1.94%  "

# Row 10
Set-TextValue "E10" "  +5.94%  "

# Row 11
Set-TextValue "D11" "7.33"
Set-TextValue "E11" "  +5.42%  "

# Row 12
Set-TextValue "E12" "  +3.58%  "

# Row 13
Set-TextValue "B13" "Avalanche"
Set-TextValue "C13" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D13" "33.21"
Set-TextValue "E13" "  +4.89%  "

# Row 14
Set-TextValue "B14" "ShibaInu"
Set-TextValue "C14" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D14" "0.0000221"
Set-TextValue "E14" "  +1.85%  "

# Row 15
Set-TextValue "D15" "4.151.93"
Set-TextValue "E15" "  +1.83%  "

# Row 16
Set-TextValue "D16" "3.549.25"
Set-TextValue "E16" "  +1.70%  "

# Row 17
Set-TextValue "D17" "68.022.94"
Set-TextValue "E17" "  +1.15%  "

# Row 18
Set-TextValue "E18" "  +0.28%  "

# Row 19
Set-TextValue "D19" "6.80"
Set-TextValue "E19" "  +6.46%  "

# Row 20
Set-TextValue "D20" "15.94"
Set-TextValue "E20" "  +5.98%  "

# Row 21
Set-TextValue "D21" "9.94"
Set-TextValue "E21" "  +10.23%  "

# Row 22
Set-TextValue "D22" "455.43"
Set-TextValue "E22" "  +1.88%  "

# Row 23
Set-TextValue "D23" "0.642"
Set-TextValue "E23" "  +3.27%  "

# Row 24
Set-TextValue "D24" "78.29"
Set-TextValue "E24" "  +1.40%  "

# Row 25
Set-TextValue "D25" "10.54"
Set-TextValue "E25" "  +3.92%  "

# Row 26
Set-TextValue "E26" "  +1.60%  "

# Row 27
Set-TextValue "D27" "3.690.94"
Set-TextValue "E27" "  +1.76%  "

# Row 28
Set-TextValue "E28" "  -0.09%  "

# Row 29
Set-TextValue "E29" "  +8.85%  "

# Row 30
Set-TextValue "D30" "2.57"
Set-TextValue "E30" "  +3.65%  "

# Row 31
Set-TextValue "D31" "1.69"
Set-TextValue "E31" "  +7.35%  "

# Row 32
Set-TextValue "D32" "0.171"
Set-TextValue "E32" "  +3.87%  "

# Row 33
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  -0.02%  "

# Row 34
Set-TextValue "D34" "6.40"
Set-TextValue "E34" "  +4.64%  "

# Row 35
Set-TextValue "D35" "26.08"
Set-TextValue "E35" "  +1.44%  "

# Row 36
Set-TextValue "E36" "  +3.99%  "

# Row 37
Set-TextValue "D37" "3.544.48"
Set-TextValue "E37" "  +1.87%  "

# Row 38
Set-TextValue "D38" "8.27"
Set-TextValue "E38" "  +3.68%  "

# Row 39
Set-TextValue "D39" "2.38"
Set-TextValue "E39" "  +7.98%  "

# Row 40
Set-TextValue "E40" "  +0.06%  "

# Row 41
Set-TextValue "D41" "178.18"
Set-TextValue "E41" "  +3.37%  "

# Row 42
Set-TextValue "D42" "0.0919"
Set-TextValue "E42" "  +5.68%  "

# Row 43
Set-TextValue "D43" "0.999"
Set-TextValue "E43" "  -0.16%  "

# Row 44
Set-TextValue "D44" "5.60"
Set-TextValue "E44" "  +3.39%  "

# Row 45
Set-TextValue "D45" "31.10"
Set-TextValue "E45" "  +16.01%  "

# Row 46
Set-TextValue "E46" "  +1.30%  "

# Row 47
Set-TextValue "B47" "OKB"
Set-TextValue "C47" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D47" "46.64"
Set-TextValue "E47" "  +2.72%  "

# Row 48
Set-TextValue "B48" "ONDO"
Set-TextValue "C48" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D48" "1.33"
Set-TextValue "E48" "  +7.08%  "

# Row 49
Set-TextValue "D49" "2.65"
Set-TextValue "E49" "  +3.98%  "

# Row 50
Set-TextValue "D50" "7.79"
Set-TextValue "E50" "  +3.39%  "

# Row 51
Set-TextValue "E51" "  +1.95%  "
